$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F57").Value = "원전공 기준의 등록금을 납부하면 됨. 다전공을 신청해도 추가 등록금을 내지 않음. 다전공 제도 참여시 별도의 장학금 지원은 현재 없음. 등록금 고지, 환불 등에 관한 사항은 학교 홈페이지 한경공지 참조. https://www.hknu.ac.kr/bbs/kor/69/87311/artclView.do?layout=unknown 또는 재무팀 031-670-5472 연락."

$ws.Range("F58").Value = "매 학년도 전기 졸업식은 2월 중순, 후기 졸업식은 8월 중순에 진행. 정확한 일정은 학교 홈페이지 학사공지 참조. https://www.hknu.ac.kr/bbs/kor/70/87507/artclView.do?layout=unknown 또는 학사지원팀 031-670-5036 연락"

$ws.Range("F63").Value = "졸업 유예와 조기 졸업은 2월 졸업은 1월에 접수, 8월 졸업은 7월에 접수. 졸업유예는 https://www.hknu.ac.kr/bbs/kor/70/91130/artclView.do?layout=unknown 참고, 조기졸업은 https://www.hknu.ac.kr/bbs/kor/70/91129/artclView.do?layout=unknown 참고. 또는 학사지원팀 031-670-5035 연락"

$ws.Range("F65").Value = "소단위전공과정(마이크로디그리)는 식품품질관리 MD, 반려동물 MD 등 37개의 과정이 있음. 과정별로 궁금하면 마이크로디그리 목록 알려줘라고 질문. 더 자세한 사항은 왼쪽 메뉴의 다전공 제도 안내 참고. 대학 홈페이지 학사공지 참고. https://www.hknu.ac.kr/bbs/kor/70/90689/artclView.do?layout=unknown"

$wb.Save()
